# Weekly update: insert the newest week's two price records (Primera /
# Segunda quality) for "Brócoli" at Vega Central Mapocho de Santiago,
# pushing the existing history down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (row 339).
$ws.Range("A339:A340").EntireRow.Insert()

# New row 339: Brócoli, Primera, Región Metropolitana
$ws.Range("A339").Value = 9
$ws.Range("B339").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C339").Value = "Metropolitana"
$ws.Range("D339").Value = 44461
$ws.Range("E339").Value = 13
$ws.Range("F339").Value = 100112023
$ws.Range("G339").Value = "Brócoli"
$ws.Range("H339").Value = "Sin especificar"
$ws.Range("I339").Value = "Primera"
$ws.Range("J339").Value = 2500
$ws.Range("K339").Value = 600
$ws.Range("L339").Value = 650
$ws.Range("M339").Value = 625
$ws.Range("N339").Value = "$/unidad"
$ws.Range("O339").Value = "Región Metropolitana"
$ws.Range("P339").Value = 625
$ws.Range("Q339").Value = 1
$ws.Range("R339").Value = "Hortaliza"

# New row 340: Brócoli, Segunda, Región Metropolitana
$ws.Range("A340").Value = 9
$ws.Range("B340").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C340").Value = "Metropolitana"
$ws.Range("D340").Value = 44461
$ws.Range("E340").Value = 13
$ws.Range("F340").Value = 100112023
$ws.Range("G340").Value = "Brócoli"
$ws.Range("H340").Value = "Sin especificar"
$ws.Range("I340").Value = "Segunda"
$ws.Range("J340").Value = 790
$ws.Range("K340").Value = 450
$ws.Range("L340").Value = 500
$ws.Range("M340").Value = 475
$ws.Range("N340").Value = "$/unidad"
$ws.Range("O340").Value = "Región Metropolitana"
$ws.Range("P340").Value = 475
$ws.Range("Q340").Value = 1
$ws.Range("R340").Value = "Hortaliza"
